$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 (header / sample-size counts)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 (CON)
$ws.Range("B2").Value = 2.3155599477277584
$ws.Range("C2").Value = 8.3012195145644281
$ws.Range("D2").Value = 9.609498529532841
$ws.Range("E2").Value = 9.013385994064695

# Update row 3 (STR)
$ws.Range("B3").Value = 2.3277528707258512
$ws.Range("C3").Value = 5.7196125950586341
$ws.Range("D3").Value = 13.251099650155393
$ws.Range("E3").Value = 5.3717862007717372

# Match the new selection recorded in the saved workbook
$ws.Range("B1:E3").Select()
